$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Floors" column (G) values: expand compact ranges like "3-11" into
# explicit comma-separated floor lists, per updated unit testing on floor number.
$ws.Range("G2").Value = "4-9"
$ws.Range("G3").Value = "9-32"
$ws.Range("G12").Value = "9-11,8,3-7"
$ws.Range("G13").Value = "9,11,13,15,17,19,21,23"
$ws.Range("G14").Value = "14,16,18,20,22"
$ws.Range("G15").Value = "2,3,6,7,10,11"
$ws.Range("G16").Value = "3-7"
$ws.Range("G17").Value = "4,6,8,10,12,14,16,18,20,22,24,26,28"
$ws.Range("G18").Value = "13,15,17,19,21,23"
$ws.Range("G19").Value = "3,5,7,9,11,13,15,17,19,21,23"
$ws.Range("G20").Value = "8"
$ws.Range("G21").Value = "10,12,14,16,18,20,22,24,26,28"
$ws.Range("G22").Value = "2,3-7,8,9-11"
$ws.Range("G23").Value = "24,26,28"
$ws.Range("G24").Value = "25,27"
$ws.Range("G25").Value = "24,26,28"

# Update selection to reflect the author's cursor position after editing.
$ws.Range("G26").Select()
